$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.397.77"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "'1.858.25"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'314.74"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "'0.4631"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.3724"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").Value = "'0.07321"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "'0.8890"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("D11").Value = "'19.98"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "'0.07845"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("D13").Value = "'1.931.29"
$ws.Range("E13").Value = "  +8.48%  "
$ws.Range("D14").Value = "'5.408"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "'6.575"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "'91.92"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "'0.000008982"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "'14.79"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "'27.409.68"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'5.133"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'10.55"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "'2.089.22"
$ws.Range("E24").Value = "  +4.51%  "
$ws.Range("D25").Value = "'1.932"
$ws.Range("E25").Value = "  +4.23%  "
$ws.Range("D26").Value = "'151.99"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").Value = "'18.45"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").Value = "'2.061"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "'5.102"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").Value = "'116.20"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "'0.08857"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'3.154"
$ws.Range("E32").Value = "  +6.17%  "
$ws.Range("D33").Value = "'0.7705"
$ws.Range("E33").Value = "  +5.06%  "
$ws.Range("D34").Value = "'1.177"
$ws.Range("E34").Value = "  +3.41%  "
$ws.Range("D35").Value = "'4.517"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").Value = "'2.689"
$ws.Range("E36").Value = "  +8.61%  "
$ws.Range("D37").Value = "'1.082"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "'0.01964"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "'0.05241"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "'2.959"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "'7.081"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "'0.5156"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'0.1640"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").Value = "'8.422"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").Value = "'0.4820"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "'10.35"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'103.00"
$ws.Range("D49").Value = "'1.655"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").Value = "'0.06220"
$ws.Range("D51").Value = "'65.26"
$ws.Range("E51").Value = "  -0.26%  "
